$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of Train Run Trends data for 5/9 (Excel serial date 42499)
$ws.Range("A26").Value = 42499
$ws.Range("B26").Value = 142
$ws.Range("C26").Value = 137
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 137
$ws.Range("G26").Value = 0.96478873239436624
$ws.Range("H26").Value = 43.825352112787222
$ws.Range("I26").Value = 36.116666675079614
$ws.Range("J26").Value = 58.833333340007812

# Match the author's final selection after entering the new data
$ws.Range("C30").Select() | Out-Null
